$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44168
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 430
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = 440
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 440

# Row 3
$ws.Range("D3").Value = 44575
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 380
$ws.Range("L3").Value = 400
$ws.Range("M3").Value = 390
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 390

# Row 4
$ws.Range("D4").Value = 44523
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 550
$ws.Range("L4").Value = 580
$ws.Range("M4").Value = 565
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 565

# Row 7
$ws.Range("D7").Value = 44251
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 250
$ws.Range("L7").Value = 280
$ws.Range("M7").Value = 265
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 265

# Row 8
$ws.Range("D8").Value = 44224
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 230
$ws.Range("L8").Value = 250
$ws.Range("M8").Value = 240
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 240

# Row 9
$ws.Range("D9").Value = 44224
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 230
$ws.Range("M9").Value = 215
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 215

# Row 10
$ws.Range("D10").Value = 44580
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 380
$ws.Range("L10").Value = 400
$ws.Range("M10").Value = 390
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 390

# Row 11
$ws.Range("D11").Value = 44253
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 270
$ws.Range("L11").Value = 280
$ws.Range("M11").Value = 275
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 275

# Row 12
$ws.Range("D12").Value = 44176
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1300
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = 375
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 375

# Row 13
$ws.Range("D13").Value = 44586
$ws.Range("I13").Value = "Tercera"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 330
$ws.Range("L13").Value = 350
$ws.Range("M13").Value = 340
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 340

# Row 14
$ws.Range("D14").Value = 44589
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 325
$ws.Range("L14").Value = 350
$ws.Range("M14").Value = 338
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 338

# Row 15
$ws.Range("D15").Value = 44202
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 1300
$ws.Range("K15").Value = 230
$ws.Range("L15").Value = 250
$ws.Range("M15").Value = 240
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 240

# Row 16
$ws.Range("D16").Value = 44172
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 420
$ws.Range("M16").Value = 410
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 410

# Row 17
$ws.Range("D17").Value = 44547
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 350
$ws.Range("L17").Value = 370
$ws.Range("M17").Value = 360
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 360

# Row 18
$ws.Range("D18").Value = 44217
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 1600
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 350
$ws.Range("M18").Value = 325
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 325

# Row 19
$ws.Range("D19").Value = 44166
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 1700
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 530
$ws.Range("M19").Value = 515
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 515

# Row 20
$ws.Range("D20").Value = 44175
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 1200
$ws.Range("K20").Value = 400
$ws.Range("L20").Value = 430
$ws.Range("M20").Value = 415
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 415

# Row 21
$ws.Range("D21").Value = 44201
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 1800
$ws.Range("K21").Value = 250
$ws.Range("L21").Value = 270
$ws.Range("M21").Value = 260
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 260

# Row 22
$ws.Range("D22").Value = 44566
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 320
$ws.Range("M22").Value = 310
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 310

# Row 23
$ws.Range("D23").Value = 44162
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 900
$ws.Range("K23").Value = 500
$ws.Range("L23").Value = 550
$ws.Range("M23").Value = 525
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 525

# Row 24
$ws.Range("D24").Value = 44162
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 1200
$ws.Range("K24").Value = 500
$ws.Range("L24").Value = 550
$ws.Range("M24").Value = 525
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 525

# Row 25
$ws.Range("D25").Value = 44231
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 180
$ws.Range("L25").Value = 200
$ws.Range("M25").Value = 190
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 190

# Row 26
$ws.Range("D26").Value = 44530
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 1300
$ws.Range("K26").Value = 450
$ws.Range("L26").Value = 480
$ws.Range("M26").Value = 465
$ws.Range("O26").Value = "Perú"
$ws.Range("P26").Value = 465

# Row 27
$ws.Range("D27").Value = 44301
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 280
$ws.Range("L27").Value = 300
$ws.Range("M27").Value = 290
$ws.Range("O27").Value = "Perú"
$ws.Range("P27").Value = 290

# Row 28
$ws.Range("D28").Value = 44229
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 230
$ws.Range("L28").Value = 250
$ws.Range("M28").Value = 240
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 240

# Row 29
$ws.Range("D29").Value = 44214
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 1200
$ws.Range("K29").Value = 400
$ws.Range("L29").Value = 450
$ws.Range("M29").Value = 425
$ws.Range("O29").Value = "Perú"
$ws.Range("P29").Value = 425

# Row 30
$ws.Range("D30").Value = 44160
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 550
$ws.Range("M30").Value = 525
$ws.Range("O30").Value = "Perú"
$ws.Range("P30").Value = 525
